$d = $word.ActiveDocument

# --- Step 1: relocate the _GoBack bookmark from "Dit is de structuur van lo|caties."
#     to inside "is|Nothing" in the hasValue-combinator paragraph, and merge the
#     two runs it used to separate back into a single run.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# the old bookmark used to split "Dit is de structuur van lo|caties." into two
# runs; with the bookmark gone, merge them back into a single run.
$r = $d.Content
$null = $r.Find.Execute("Dit is de structuur van locaties.", $false, $false, $false, $false, $false, $true, 1, $false, "Dit is de structuur van locaties.", 2)

$r = $d.Content
$null = $r.Find.Execute("isNothing singleGetter. De definitie luidt:", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$bookmarkPos = $r.Start + 2
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 2: ". De definitie luidt:" -> ". De definitie is in essentie:"
$r = $d.Content
$null = $r.Find.Execute(". De definitie luidt:", $true, $false, $false, $false, $false, $true, 1, $false, ". De definitie is in essentie:", 2)

# --- Step 3: "f >-> isNothing"  -> "f >=> isNothing"  (only touch the operator)
$r = $d.Content
$null = $r.Find.Execute("f >-> isNothing", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$opRange = $d.Range($r.Start + 1, $r.End - 9)
$opRange.Text = " >=> "

# --- Step 4: prepend the new parenthetical sentence in front of
#     "De structuur van locaties die gegenereerd wordt, is dat van single to single."
$r = $d.Content
$null = $r.Find.Execute("De structuur van locaties die gegenereerd wordt, is dat van single", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$insPoint = $d.Range($r.Start, $r.Start)
$insPoint.InsertBefore("(De werkelijke definitie is iets ingewikkelder omdat isNothing geen SingleGetter is.) ")

# style "isNothing" (within the sentence just inserted) as inline code
$r = $d.Content
$null = $r.Find.Execute("omdat isNothing geen", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$codeStart = $r.Start + 6
$codeRange = $d.Range($codeStart, $codeStart + 9)
$codeRange.Style = "inlinecode"

# --- Step 5: " gegenereerd wordt, is dat van single to single." ->
#     " gegenereerd wordt, is dus die van kleisli composition, oftewel die van LocatonT bind."
$r = $d.Content
$null = $r.Find.Execute(" gegenereerd wordt, is dat van single to single.", $false, $false, $false, $false, $false, $true, 1, $false, " gegenereerd wordt, is dus die van kleisli composition, oftewel die van LocatonT bind.", 2)

# style "LocatonT bind" as inline code
$r = $d.Content
$null = $r.Find.Execute("oftewel die van LocatonT bind.", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$codeStart2 = $r.Start + 16
$codeRange2 = $d.Range($codeStart2, $codeStart2 + 13)
$codeRange2.Style = "inlinecode"

# --- Step 6: footer page field cached value 8 -> 6
$f = $d.Sections(1).Footers(1).Range
$null = $f.Find.Execute("8 van 8", $true, $false, $false, $false, $false, $true, 1, $false, "6 van 8", 2)
